# Updates the cryptos list with refreshed prices / 1h volume change percentages
# (and restores the Wrapped Ether / Wrapped liquid staked Ether 2.0 row order).
#
# D/E columns hold plain text (e.g. "25.768.57", "  -0.09%  ") rather than numbers,
# so each cell is forced to Text format before the value is written and the format
# is reset back to the sheet default ("Normal" style) afterwards -- this stops Excel
# from auto-converting decimal-looking strings into real numbers while keeping the
# cell style identical to the surrounding (unstyled) data cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '25.768.57'
Set-TextValue "E2" '  -0.09%  '

# Row 3
Set-TextValue "D3" '1.636.43'
Set-TextValue "E3" '  -0.09%  '

# Row 4
Set-TextValue "E4" '  +0.13%  '

# Row 5
Set-TextValue "D5" '215.64'
Set-TextValue "E5" '  +0.06%  '

# Row 6
Set-TextValue "D6" '0.502'
Set-TextValue "E6" '  -0.78%  '

# Row 7
Set-TextValue "E7" '  +0.15%  '

# Row 8
Set-TextValue "E8" '  +0.09%  '

# Row 9
Set-TextValue "D9" '0.0637'
Set-TextValue "E9" '  -0.86%  '

# Row 10
Set-TextValue "D10" '19.59'
Set-TextValue "E10" '  -3.82%  '

# Row 11
Set-TextValue "E11" '  +0.99%  '

# Row 12
Set-TextValue "D12" '4.25'
Set-TextValue "E12" '  -0.39%  '

# Row 13
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.646.74'
Set-TextValue "E13" '  +0.53%  '

# Row 14
Set-TextValue "B14" 'WrappedliquidstakedEther2.0'
Set-TextValue "C14" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D14" '1.861.36'
Set-TextValue "E14" '  -0.11%  '

# Row 15
Set-TextValue "D15" '0.555'
Set-TextValue "E15" '  -1.13%  '

# Row 16
Set-TextValue "D16" '0.0₃0767'
Set-TextValue "E16" '  +0.33%  '

# Row 17
Set-TextValue "D17" '62.85'
Set-TextValue "E17" '  -0.69%  '

# Row 18
Set-TextValue "D18" '25.782.07'
Set-TextValue "E18" '  -0.11%  '

# Row 19
Set-TextValue "E19" '  +0.06%  '

# Row 21
Set-TextValue "D21" '194.10'
Set-TextValue "E21" '  +0.54%  '

# Row 22
Set-TextValue "D22" '9.96'
Set-TextValue "E22" '  +0.47%  '

# Row 23
Set-TextValue "D23" '6.28'
Set-TextValue "E23" '  +2.45%  '

# Row 24
Set-TextValue "E24" '  +0.12%  '

# Row 25
Set-TextValue "D25" '1.83'
Set-TextValue "E25" '  +1.96%  '

# Row 26
Set-TextValue "D26" '140.05'

# Row 27
Set-TextValue "E27" '  -2.12%  '

# Row 28
Set-TextValue "D28" '6.87'

# Row 29
Set-TextValue "D29" '15.54'
Set-TextValue "E29" '  -0.17%  '

# Row 30
Set-TextValue "E30" '  -0.09%  '

# Row 31
Set-TextValue "D31" '0.0495'
Set-TextValue "E31" '  -0.16%  '

# Row 32
Set-TextValue "D32" '3.33'
Set-TextValue "E32" '  +1.38%  '

# Row 33
Set-TextValue "D33" '3.25'
Set-TextValue "E33" '  +0.42%  '

# Row 34
Set-TextValue "E34" '  +1.49%  '

# Row 35
Set-TextValue "E35" '  +0.32%  '

# Row 36
Set-TextValue "D36" '0.900'
Set-TextValue "E36" '  -0.37%  '

# Row 37
Set-TextValue "D37" '0.549'
Set-TextValue "E37" '  -1.64%  '

# Row 38
Set-TextValue "D38" '1.116.84'
Set-TextValue "E38" '  -1.28%  '

# Row 39
Set-TextValue "E39" '  -2.16%  '

# Row 40
Set-TextValue "E40" '  -0.40%  '

# Row 41
Set-TextValue "E41" '  +0.70%  '

# Row 42
Set-TextValue "E42" '  +1.46%  '

# Row 43
Set-TextValue "D43" '99.73'
Set-TextValue "E43" '  +0.86%  '

# Row 44
Set-TextValue "D44" '0.802'
Set-TextValue "E44" '  +0.01%  '

# Row 45
Set-TextValue "D45" '1.771.11'
Set-TextValue "E45" '  -0.21%  '

# Row 46
Set-TextValue "D46" '0.0₆0108'
Set-TextValue "E46" '  -1.76%  '

# Row 47
Set-TextValue "D47" '55.21'
Set-TextValue "E47" '  -0.74%  '

# Row 48
Set-TextValue "D48" '0.417'
Set-TextValue "E48" '  -2.23%  '

# Row 49
Set-TextValue "D49" '0.0501'
Set-TextValue "E49" '  -0.36%  '

# Row 50
Set-TextValue "E50" '  -1.90%  '

# Row 51
Set-TextValue "E51" '  +2.73%  '
